$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.574.23"
$ws.Range("E2").Value = "  +1.12%  "

$ws.Range("D3").Value = "1.855.13"

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.86"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4735"
$ws.Range("E7").Value = "  +0.68%  "

$ws.Range("E8").Value = "  +1.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06322"
$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.70"
$ws.Range("E10").Value = "  +8.78%  "

$ws.Range("D11").Value = "1.881.16"
$ws.Range("E11").Value = "  +1.63%  "

$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.000"
$ws.Range("E13").Value = "  +1.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.58"
$ws.Range("E14").Value = "  -0.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6260"
$ws.Range("E15").Value = "  -0.17%  "

$ws.Range("D16").Value = "30.536.88"
$ws.Range("E16").Value = "  +1.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "245.36"
$ws.Range("E17").Value = "  +7.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9995"
$ws.Range("E18").Value = "  -0.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.71"
$ws.Range("E19").Value = "  +0.94%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007336"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9991"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.941"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.946"
$ws.Range("E23").Value = "  +0.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.142"
$ws.Range("E24").Value = "  -0.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.88"
$ws.Range("E25").Value = "  -2.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.03"
$ws.Range("E26").Value = "  +1.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.880"
$ws.Range("E27").Value = "  +0.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1018"
$ws.Range("E28").Value = "  -1.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.358"
$ws.Range("E29").Value = "  -1.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.014"
$ws.Range("E30").Value = "  -2.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.837"
$ws.Range("E31").Value = "  -0.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04841"
$ws.Range("E32").Value = "  -0.93%  "

$ws.Range("E33").Value = "  -1.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7026"
$ws.Range("E34").Value = "  -0.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.691"
$ws.Range("E35").Value = "  -0.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01898"
$ws.Range("E36").Value = "  +1.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.677"
$ws.Range("E37").Value = "  +1.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.001"
$ws.Range("E38").Value = "  +2.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8765"
$ws.Range("E39").Value = "  -3.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "106.85"
$ws.Range("E40").Value = "  +1.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.550"
$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4061"
$ws.Range("E43").Value = "  -0.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.195"
$ws.Range("E44").Value = "  +1.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.81"
$ws.Range("E45").Value = "  +4.13%  "

$ws.Range("E46").Value = "  +1.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.61"
$ws.Range("E47").Value = "  +1.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.570"
$ws.Range("E48").Value = "  -0.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05535"
$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("E50").Value = "  -2.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3683"
$ws.Range("E51").Value = "  +0.36%  "
